# Update column G ("K") values on Sheet1 with newly regenerated strikeout counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$gValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 0
    6  = 1
    7  = 2
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 3
    13 = 0
    14 = 0
    15 = 3
    16 = 0
    17 = 2
    18 = 0
    19 = 0
    20 = 1
    21 = 2
    22 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
